# Auto-generated edit script applying the scrape update (13:05:03 -> 13:53:00)
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:53:00"
$ws.Cells.Item(3, 1).Value = "Total filas: 260"
$ws.Cells.Item(43, 3).Value = "215A_EL PATO"
$ws.Cells.Item(44, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(49, 3).Value = "10_OLMOS"
$ws.Cells.Item(50, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(55, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(56, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(80, 1).Value = "07:50:28"
$ws.Cells.Item(80, 3).Value = "215B_EL PATO"
$ws.Cells.Item(80, 4).Value = 33
$ws.Cells.Item(81, 1).Value = "07:19:11"
$ws.Cells.Item(81, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(81, 4).Value = 64
$ws.Cells.Item(105, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(106, 3).Value = "17_ROMERO"
$ws.Cells.Item(107, 1).Value = "08:16:46"
$ws.Cells.Item(107, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(107, 4).Value = 67
$ws.Cells.Item(108, 1).Value = "08:51:14"
$ws.Cells.Item(108, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(108, 4).Value = 32
$ws.Cells.Item(109, 1).Value = "08:38:01"
$ws.Cells.Item(109, 3).Value = "17_ROMERO"
$ws.Cells.Item(109, 4).Value = 45
$ws.Cells.Item(131, 1).Value = "08:51:14"
$ws.Cells.Item(131, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(131, 4).Value = 80
$ws.Cells.Item(132, 1).Value = "09:28:57"
$ws.Cells.Item(132, 3).Value = "15_ABASTO"
$ws.Cells.Item(132, 4).Value = 43
$ws.Cells.Item(133, 1).Value = "08:51:14"
$ws.Cells.Item(133, 3).Value = "15_ABASTO"
$ws.Cells.Item(133, 4).Value = 81
$ws.Cells.Item(134, 1).Value = "09:28:57"
$ws.Cells.Item(134, 3).Value = "10_OLMOS"
$ws.Cells.Item(134, 4).Value = 44
$ws.Cells.Item(167, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(168, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(170, 1).Value = "11:02:02"
$ws.Cells.Item(170, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(170, 4).Value = 33
$ws.Cells.Item(171, 1).Value = "10:28:51"
$ws.Cells.Item(171, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(171, 4).Value = 67
$ws.Cells.Item(195, 3).Value = "14_ABASTO"
$ws.Cells.Item(196, 1).Value = "11:44:55"
$ws.Cells.Item(196, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(196, 4).Value = 36
$ws.Cells.Item(197, 1).Value = "12:04:34"
$ws.Cells.Item(197, 3).Value = "215A_EL PATO"
$ws.Cells.Item(197, 4).Value = 16
$ws.Cells.Item(201, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(202, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(225, 1).Value = "12:04:34"
$ws.Cells.Item(225, 3).Value = "14_ABASTO"
$ws.Cells.Item(225, 4).Value = 82
$ws.Cells.Item(226, 1).Value = "13:05:03"
$ws.Cells.Item(226, 3).Value = "15_ABASTO"
$ws.Cells.Item(226, 4).Value = 21
$ws.Cells.Item(237, 1).Value = "13:53:00"
$ws.Cells.Item(237, 2).Value = "13:55"
$ws.Cells.Item(237, 4).Value = 2
$ws.Cells.Item(238, 1).Value = "13:05:03"
$ws.Cells.Item(238, 3).Value = "225_GOMEZ"
$ws.Cells.Item(238, 4).Value = 51
$ws.Cells.Item(239, 1).Value = "12:04:34"
$ws.Cells.Item(239, 2).Value = "13:56"
$ws.Cells.Item(239, 4).Value = 112
$ws.Cells.Item(240, 1).Value = "13:53:00"
$ws.Cells.Item(240, 2).Value = "13:57"
$ws.Cells.Item(240, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(240, 4).Value = 4
$ws.Cells.Item(241, 1).Value = "13:53:00"
$ws.Cells.Item(241, 2).Value = "14:04"
$ws.Cells.Item(241, 3).Value = "17_ROMERO"
$ws.Cells.Item(241, 4).Value = 11
$ws.Cells.Item(242, 1).Value = "13:53:00"
$ws.Cells.Item(242, 2).Value = "14:05"
$ws.Cells.Item(242, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(242, 4).Value = 12
$ws.Cells.Item(243, 1).Value = "13:53:00"
$ws.Cells.Item(243, 2).Value = "14:07"
$ws.Cells.Item(243, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(243, 4).Value = 14
$ws.Cells.Item(244, 1).Value = "13:53:00"
$ws.Cells.Item(244, 2).Value = "14:09"
$ws.Cells.Item(244, 3).Value = "10_OLMOS"
$ws.Cells.Item(244, 4).Value = 16
$ws.Cells.Item(245, 1).Value = "13:53:00"
$ws.Cells.Item(245, 2).Value = "14:12"
$ws.Cells.Item(245, 3).Value = "15_ABASTO"
$ws.Cells.Item(245, 4).Value = 19
$ws.Cells.Item(246, 1).Value = "13:53:00"
$ws.Cells.Item(246, 2).Value = "14:17"
$ws.Cells.Item(246, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(246, 4).Value = 24
$ws.Cells.Item(247, 1).Value = "13:53:00"
$ws.Cells.Item(247, 2).Value = "14:20"
$ws.Cells.Item(247, 3).Value = "215C_EL PATO"
$ws.Cells.Item(247, 4).Value = 27
$ws.Cells.Item(248, 1).Value = "13:53:00"
$ws.Cells.Item(248, 2).Value = "14:21"
$ws.Cells.Item(248, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(248, 4).Value = 28
$ws.Cells.Item(249, 1).Value = "13:53:00"
$ws.Cells.Item(249, 2).Value = "14:27"
$ws.Cells.Item(249, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(249, 4).Value = 34
$ws.Cells.Item(250, 1).Value = "13:53:00"
$ws.Cells.Item(250, 2).Value = "14:31"
$ws.Cells.Item(250, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(250, 4).Value = 38
$ws.Cells.Item(250, 5).Value = "LP1912"
$ws.Cells.Item(251, 1).Value = "13:53:00"
$ws.Cells.Item(251, 2).Value = "14:45"
$ws.Cells.Item(251, 3).Value = "14_ABASTO"
$ws.Cells.Item(251, 4).Value = 52
$ws.Cells.Item(251, 5).Value = "LP1912"
$ws.Cells.Item(252, 1).Value = "13:05:03"
$ws.Cells.Item(252, 2).Value = "14:48"
$ws.Cells.Item(252, 3).Value = "14_ABASTO"
$ws.Cells.Item(252, 4).Value = 103
$ws.Cells.Item(252, 5).Value = "LP1912"
$ws.Cells.Item(253, 1).Value = "13:53:00"
$ws.Cells.Item(253, 2).Value = "14:57"
$ws.Cells.Item(253, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(253, 4).Value = 64
$ws.Cells.Item(253, 5).Value = "LP1912"
$ws.Cells.Item(254, 1).Value = "13:53:00"
$ws.Cells.Item(254, 2).Value = "14:58"
$ws.Cells.Item(254, 3).Value = "215B_EL PATO"
$ws.Cells.Item(254, 4).Value = 65
$ws.Cells.Item(254, 5).Value = "LP1912"
$ws.Cells.Item(255, 1).Value = "13:53:00"
$ws.Cells.Item(255, 2).Value = "15:00"
$ws.Cells.Item(255, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(255, 4).Value = 67
$ws.Cells.Item(255, 5).Value = "LP1912"
$ws.Cells.Item(256, 1).Value = "13:53:00"
$ws.Cells.Item(256, 2).Value = "15:05"
$ws.Cells.Item(256, 3).Value = "10_OLMOS"
$ws.Cells.Item(256, 4).Value = 72
$ws.Cells.Item(256, 5).Value = "LP1912"
$ws.Cells.Item(257, 1).Value = "13:53:00"
$ws.Cells.Item(257, 2).Value = "15:10"
$ws.Cells.Item(257, 3).Value = "17_ROMERO"
$ws.Cells.Item(257, 4).Value = 77
$ws.Cells.Item(257, 5).Value = "LP1912"
$ws.Cells.Item(258, 1).Value = "13:53:00"
$ws.Cells.Item(258, 2).Value = "15:14"
$ws.Cells.Item(258, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(258, 4).Value = 81
$ws.Cells.Item(258, 5).Value = "LP1912"
$ws.Cells.Item(259, 1).Value = "13:53:00"
$ws.Cells.Item(259, 2).Value = "15:20"
$ws.Cells.Item(259, 3).Value = "15_ABASTO"
$ws.Cells.Item(259, 4).Value = 87
$ws.Cells.Item(259, 5).Value = "LP1912"
$ws.Cells.Item(260, 1).Value = "13:53:00"
$ws.Cells.Item(260, 2).Value = "15:23"
$ws.Cells.Item(260, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(260, 4).Value = 90
$ws.Cells.Item(260, 5).Value = "LP1912"
$ws.Cells.Item(261, 1).Value = "13:53:00"
$ws.Cells.Item(261, 2).Value = "15:32"
$ws.Cells.Item(261, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(261, 4).Value = 99
$ws.Cells.Item(261, 5).Value = "LP1912"
$ws.Cells.Item(262, 1).Value = "13:53:00"
$ws.Cells.Item(262, 2).Value = "15:35"
$ws.Cells.Item(262, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(262, 4).Value = 102
$ws.Cells.Item(262, 5).Value = "LP1912"
$ws.Cells.Item(263, 1).Value = "13:53:00"
$ws.Cells.Item(263, 2).Value = "15:39"
$ws.Cells.Item(263, 3).Value = "215A_EL PATO"
$ws.Cells.Item(263, 4).Value = 106
$ws.Cells.Item(263, 5).Value = "LP1912"
$ws.Cells.Item(264, 1).Value = "13:53:00"
$ws.Cells.Item(264, 2).Value = "15:47"
$ws.Cells.Item(264, 3).Value = "14_ABASTO"
$ws.Cells.Item(264, 4).Value = 114
$ws.Cells.Item(264, 5).Value = "LP1912"
$ws.Cells.Item(265, 1).Value = "13:53:00"
$ws.Cells.Item(265, 2).Value = "15:47"
$ws.Cells.Item(265, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(265, 4).Value = 114
$ws.Cells.Item(265, 5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:53:00"
$ws.Cells.Item(3, 1).Value = "Total filas: 33"
$ws.Cells.Item(36, 1).Value = "13:53:00"
$ws.Cells.Item(36, 4).Value = 27
$ws.Cells.Item(37, 1).Value = "13:53:00"
$ws.Cells.Item(37, 4).Value = 65
$ws.Cells.Item(38, 1).Value = "13:53:00"
$ws.Cells.Item(38, 2).Value = "15:39"
$ws.Cells.Item(38, 3).Value = "215A_EL PATO"
$ws.Cells.Item(38, 4).Value = 106
$ws.Cells.Item(38, 5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:53:00"
$ws.Cells.Item(3, 1).Value = "Total filas: 30"
$ws.Cells.Item(33, 1).Value = "13:53:00"
$ws.Cells.Item(33, 4).Value = 16
$ws.Cells.Item(34, 1).Value = "13:53:00"
$ws.Cells.Item(34, 4).Value = 60
$ws.Cells.Item(35, 1).Value = "13:53:00"
$ws.Cells.Item(35, 2).Value = "15:34"
$ws.Cells.Item(35, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(35, 4).Value = 101
$ws.Cells.Item(35, 5).Value = "L6173"

Write-Host "Edit script completed successfully."
